# Regenerate the s_vals data to filter save games.
# Updates numeric columns B:G for rows 2-8 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: B, C, D, E, F  (G is the row sum: B+C+D+E+F)
$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0)
    4 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0)
    5 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0)
    6 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0)
    7 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0)
    8 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1)
}

$gSums = @{
    2 = 5.582307763322248
    3 = 5.582307763322248
    4 = 4.327115817150455
    5 = 6.15379541431027
    6 = 8.656069925401464
    7 = 6.15379541431027
    8 = 5.582307763322248
}

foreach ($row in 2..8) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $gSums[$row]
}
